$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mutualista Imbabura" entity was renamed/corrected to the real
# institution name in the two credit rows (U2, U3).
$newName = "Cooperativa de Ahorro y Crédito San José de Vittoria"

$rng = $ws.Range("U2:U3")
$rng.Value = $newName
$rng.Font.Name = "Calibri"
$rng.Font.Size = 11
$rng.Font.ColorIndex = 1
$rng.WrapText = $true

# Leave the selection where the author's last click landed.
$ws.Range("V7").Select()
